$wb = $excel.ActiveWorkbook

# --- Overview sheet (sheet1) : add "Info :" header block with Version/Revision/GenDate rows ---
$ws1 = $wb.Worksheets.Item("Overview")

# Order matters: it drives the order new entries are appended to the shared-string table.
$ws1.Range("A4").Value = "Version"
$ws1.Range("A3").Value = "Revision"
$ws1.Range("A5").Value = "GenDate"
$ws1.Range("A1").Value = "Info :"
$ws1.Range("A2").Value = "PN"

# Bold the info rows (A2:A5) -> creates the new bold font + cellXf
$ws1.Range("A2:A5").Font.Bold = $true

# Widen column B so later-added values have room
$ws1.Columns.Item(2).ColumnWidth = 27.8

# Make Overview the selected/active tab (this also clears tabSelected on the
# sheet that was previously active, i.e. Tasks)
$ws1.Select()
